# Auto-generated script: applies the 2023-08-02 data update
# to the violent-crime-full-year workbook (148 cell updates across 40 sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 4390
$ws.Range("J3").Value = 4655
$ws.Range("I4").Value = 1768
$ws.Range("J4").Value = 1036
$ws.Range("J5").Value = 372
$ws.Range("J6").Value = 5557
$ws.Range("I7").Value = 26215
$ws.Range("J7").Value = 16010

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 472
$ws.Range("J8").Value = 1029
$ws.Range("J10").Value = 94
$ws.Range("J11").Value = 237
$ws.Range("J15").Value = 179
$ws.Range("J16").Value = 58
$ws.Range("J20").Value = 333
$ws.Range("J23").Value = 149
$ws.Range("J26").Value = 32
$ws.Range("J29").Value = 907
$ws.Range("J31").Value = 129
$ws.Range("J33").Value = 729
$ws.Range("J37").Value = 515
$ws.Range("J40").Value = 35
$ws.Range("J41").Value = 96
$ws.Range("J42").Value = 619
$ws.Range("J47").Value = 119
$ws.Range("J49").Value = 112
$ws.Range("J51").Value = 209
$ws.Range("J53").Value = 182
$ws.Range("J54").Value = 313
$ws.Range("J63").Value = 71
$ws.Range("J65").Value = 419
$ws.Range("J66").Value = 49
$ws.Range("J67").Value = 621
$ws.Range("J74").Value = 20
$ws.Range("J79").Value = 460
$ws.Range("J83").Value = 355
$ws.Range("J85").Value = 716
$ws.Range("J86").Value = 99
$ws.Range("J88").Value = 175
$ws.Range("J89").Value = 214
$ws.Range("I91").Value = 277
$ws.Range("J91").Value = 178
$ws.Range("J92").Value = 50
$ws.Range("J95").Value = 249
$ws.Range("J97").Value = 123
$ws.Range("J99").Value = 238
$ws.Range("I101").Value = 26215
$ws.Range("J101").Value = 16010

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 244
$ws.Range("J5").Value = 31
$ws.Range("J6").Value = 240
$ws.Range("J7").Value = 729

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 94

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 276
$ws.Range("J3").Value = 314
$ws.Range("J4").Value = 52
$ws.Range("J6").Value = 226
$ws.Range("J7").Value = 907

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 181
$ws.Range("J3").Value = 264
$ws.Range("J4").Value = 56
$ws.Range("J7").Value = 716

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 90
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 130
$ws.Range("J6").Value = 309
$ws.Range("J7").Value = 619

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 150
$ws.Range("J3").Value = 183
$ws.Range("J5").Value = 20
$ws.Range("J7").Value = 515

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 82
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 147
$ws.Range("J3").Value = 247
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 621

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 103
$ws.Range("J6").Value = 181

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 289
$ws.Range("J3").Value = 317
$ws.Range("J6").Value = 336
$ws.Range("J7").Value = 1029

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 122
$ws.Range("J5").Value = 11
$ws.Range("J7").Value = 419

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 460

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 141
$ws.Range("J5").Value = 12
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 472

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 98
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 333

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 57
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 277
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 54
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 46
$ws.Range("J6").Value = 32
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 77
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 313

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 134
$ws.Range("J7").Value = 355

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 56
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 97
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 69
$ws.Range("J7").Value = 214

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 35
